# Weekly driver report update for 2025-05-05
# Updates the "Good Drivers" table (rows 13-24) on the "Driver Summary" sheet:
#   - refreshes the numeric roaming counters (client/critical/warning/good counts
#     and totals) for each adapter-driver row
#   - a few rows (14-16 and 18-21) were re-sorted by driver version, so the
#     adapter-driver label (columns A and H) and vintage date (column J) for
#     those rows move along with their new numbers
#
# Column J holds vintage dates written as plain text (e.g. "2023-03-28").
# Excel's Range.Value setter auto-detects such strings as dates and
# reformats the cell, so text writes go through Set-TextValue below, which
# temporarily forces a text number format (preventing the date inference)
# and then restores the cell's original style afterward.

function Set-TextValue($range, [string]$text) {
    $savedStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = $savedStyle
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: Intel(R) Wi-Fi 6E AX211 160MHz - 22.250.0.4 (unchanged identity/order)
$ws.Range("B13").Value = 1293197
$ws.Range("C13").Value = 4322
$ws.Range("E13").Value = 1990
$ws.Range("F13").Value = 1298528

# Row 14: now Intel(R) Wi-Fi 6E AX211 160MHz - 22.220.0.4
Set-TextValue $ws.Range("A14") "Intel(R) Wi-Fi 6E AX211 160MHz - 22.220.0.4"
$ws.Range("B14").Value = 31517
$ws.Range("C14").Value = 112
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 66
$ws.Range("F14").Value = 31629
Set-TextValue $ws.Range("H14") "22.220.0.4"
Set-TextValue $ws.Range("J14") "2023-03-28"

# Row 15: now Intel(R) Wi-Fi 6E AX211 160MHz - 23.10.0.8
Set-TextValue $ws.Range("A15") "Intel(R) Wi-Fi 6E AX211 160MHz - 23.10.0.8"
$ws.Range("B15").Value = 467311
$ws.Range("C15").Value = 772
$ws.Range("D15").Value = 931
$ws.Range("E15").Value = 706
$ws.Range("F15").Value = 469014
Set-TextValue $ws.Range("H15") "23.10.0.8"
Set-TextValue $ws.Range("J15") "2023-10-30"

# Row 16: now Intel(R) Wi-Fi 6E AX211 160MHz - 23.120.0.3
Set-TextValue $ws.Range("A16") "Intel(R) Wi-Fi 6E AX211 160MHz - 23.120.0.3"
$ws.Range("B16").Value = 455081
$ws.Range("C16").Value = 1861
$ws.Range("D16").Value = 52
$ws.Range("E16").Value = 639
$ws.Range("F16").Value = 456994
Set-TextValue $ws.Range("H16") "23.120.0.3"
Set-TextValue $ws.Range("J16") "2025-02-05"

# Row 17: Intel(R) Wi-Fi 6E AX211 160MHz - 22.230.0.8 (unchanged identity/order)
$ws.Range("B17").Value = 1787924
$ws.Range("C17").Value = 3326
$ws.Range("E17").Value = 3038
$ws.Range("F17").Value = 1793864

# Row 18: now Intel(R) Wi-Fi 6E AX211 160MHz - 23.70.2.3
Set-TextValue $ws.Range("A18") "Intel(R) Wi-Fi 6E AX211 160MHz - 23.70.2.3"
$ws.Range("B18").Value = 218767
$ws.Range("C18").Value = 334
$ws.Range("D18").Value = 313
$ws.Range("E18").Value = 573
$ws.Range("F18").Value = 219414
Set-TextValue $ws.Range("H18") "23.70.2.3"
Set-TextValue $ws.Range("J18") "2024-07-23"

# Row 19: now Intel(R) Wi-Fi 6E AX211 160MHz - 22.110.1.1
Set-TextValue $ws.Range("A19") "Intel(R) Wi-Fi 6E AX211 160MHz - 22.110.1.1"
$ws.Range("B19").Value = 135467
$ws.Range("C19").Value = 189
$ws.Range("D19").Value = 263
$ws.Range("E19").Value = 196
$ws.Range("F19").Value = 135919
Set-TextValue $ws.Range("H19") "22.110.1.1"
Set-TextValue $ws.Range("J19") "2022-01-01"

# Row 20: now Intel(R) Wi-Fi 6E AX211 160MHz - 23.100.0.4
Set-TextValue $ws.Range("A20") "Intel(R) Wi-Fi 6E AX211 160MHz - 23.100.0.4"
$ws.Range("B20").Value = 240434
$ws.Range("C20").Value = 421
$ws.Range("D20").Value = 37
$ws.Range("E20").Value = 409
$ws.Range("F20").Value = 240892
Set-TextValue $ws.Range("H20") "23.100.0.4"
Set-TextValue $ws.Range("J20") "2024-11-10"

# Row 21: now Intel(R) Wi-Fi 6E AX211 160MHz - 23.80.1.3
Set-TextValue $ws.Range("A21") "Intel(R) Wi-Fi 6E AX211 160MHz - 23.80.1.3"
$ws.Range("B21").Value = 151287
$ws.Range("C21").Value = 285
$ws.Range("D21").Value = 75
$ws.Range("E21").Value = 332
$ws.Range("F21").Value = 151647
Set-TextValue $ws.Range("H21") "23.80.1.3"
Set-TextValue $ws.Range("J21") "2024-09-03"

# Row 22: Intel(R) Wi-Fi 6E AX211 160MHz - 22.100.1.1 (unchanged identity/order)
$ws.Range("B22").Value = 272039
$ws.Range("C22").Value = 213
$ws.Range("E22").Value = 316
$ws.Range("F22").Value = 272383

# Row 23: Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.0.3 (unchanged identity/order)
$ws.Range("B23").Value = 14561
$ws.Range("E23").Value = 59
$ws.Range("F23").Value = 14561

# Row 24: Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.3.1 (unchanged identity/order)
$ws.Range("B24").Value = 12018
$ws.Range("E24").Value = 61
$ws.Range("F24").Value = 12018
